$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits after the
#    "Ath. Lituðu örvarnar ..." paragraph. Word keeps this as a hidden
#    bookmark, so it is reachable via Exists()/Item() even though it is not
#    enumerated by Bookmarks.Count.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fill in the use-case table (second table in the document).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(2)

# Row 2
$t.Cell(2, 1).Range.Text = "Bæta við nýjum notanda"
$t.Cell(2, 2).Range.Text = "Já (Verkefni 3)"

# Row 3, second cell first (simple case)
$t.Cell(3, 2).Range.Text = "Já (Verkefni 2)"

# Row 3, first cell: "Bæta við frétt/Viðburði" followed by a fresh "_GoBack"
# bookmark placed right after the text (zero-length). The engine mishandles
# zero-length bookmarks placed exactly at a paragraph's trailing edge, so we
# type a placeholder character, anchor the bookmark just before it, and then
# remove the placeholder.
$t.Cell(3, 1).Range.Text = "Bæta við frétt/ViðburðiX"
$cell31 = $t.Cell(3, 1)
$bmPos = $cell31.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$cell31b = $t.Cell(3, 1)
$delPos = $cell31b.Range.End - 2
$delRange = $d.Range($delPos, $delPos + 1)
$delRange.Text = ""

# ---------------------------------------------------------------------------
# 3) Append four new rows describing the remaining use cases.
# ---------------------------------------------------------------------------
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null

$t.Cell(4, 1).Range.Text = "Innskráning"
$t.Cell(4, 2).Range.Text = "Já (Verkefni 2)"

$t.Cell(5, 1).Range.Text = "Skoða viðburð"
$t.Cell(5, 2).Range.Text = "Já (Verkefni 2)"

$t.Cell(6, 1).Range.Text = "Skrá sig í viðburð"
$t.Cell(6, 2).Range.Text = "Já (Verkefni 2)"

$t.Cell(7, 1).Range.Text = "Útskráning"
$t.Cell(7, 2).Range.Text = "Já (Verkefni 2)"
